$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column: price text values that look numeric need an explicit Text
# number format before the assignment, otherwise Excel auto-converts the
# string to a real number (dropping trailing zeros / exact formatting).
# E column: percentage strings (with leading/trailing spaces and a % sign)
# are never auto-coerced to numbers, so a plain assignment is enough.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.899.92'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.073.64'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.16'
$ws.Range("E5").Value = '  -0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.25'
$ws.Range("E6").Value = '  -3.03%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.068.73'
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.37'
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("E11").Value = '  -1.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.470'
$ws.Range("E12").Value = '  -1.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000241'
$ws.Range("E13").Value = '  -1.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.92'
$ws.Range("E14").Value = '  -3.55%  '
$ws.Range("E15").Value = '  -1.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.586.73'
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.837.02'
$ws.Range("E17").Value = '  +0.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.00'
$ws.Range("E18").Value = '  -1.59%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.67'
$ws.Range("E19").Value = '  +2.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.073.14'
$ws.Range("E20").Value = '  -1.14%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '482.70'
$ws.Range("E21").Value = '  +1.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.687'
$ws.Range("E22").Value = '  -3.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.65'
$ws.Range("E23").Value = '  -5.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.53'
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.79'
$ws.Range("E25").Value = '  -5.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.21'
$ws.Range("E26").Value = '  -3.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.10'
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.85'
$ws.Range("E29").Value = '  -0.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.28'
$ws.Range("E30").Value = '  -6.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.61'
$ws.Range("E31").Value = '  -1.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.68'
$ws.Range("E32").Value = '  -3.25%  '
$ws.Range("E33").Value = '  -2.42%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0₃0900'
$ws.Range("E34").Value = '  -3.45%  '
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.66'
$ws.Range("E36").Value = '  -3.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.952'
$ws.Range("E37").Value = '  -2.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '45.88'
$ws.Range("E38").Value = '  -3.76%  '
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.98'
$ws.Range("E40").Value = '  -4.66%  '
$ws.Range("E41").Value = '  -2.33%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.30'
$ws.Range("E42").Value = '  -3.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.774.09'
$ws.Range("E43").Value = '  -0.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '369.41'
$ws.Range("E44").Value = '  -2.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '135.77'
$ws.Range("E45").Value = '  -0.24%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0343'
$ws.Range("E46").Value = '  -3.05%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.45'
$ws.Range("E47").Value = '  -3.01%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.41'
$ws.Range("E49").Value = '  -1.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.15'
$ws.Range("E50").Value = '  -2.51%  '
$ws.Range("E51").Value = '  -1.76%  '
